$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value2 = 3139.8823
$ws.Range("I6").Value2 = 3048.625
$ws.Range("J6").Value2 = 4600
$ws.Range("K6").Value2 = 9145.875
$ws.Range("L6").Value2 = 13800
$ws.Range("M6").Value2 = -9033.875
$ws.Range("N6").Value2 = -14024
$ws.Range("H28").Value2 = 2093.7144
$ws.Range("J28").Value2 = 1223
$ws.Range("L28").Value2 = 1223
$ws.Range("N28").Value2 = -2193
$ws.Range("H33").Value2 = 154.09091
$ws.Range("I33").Value2 = 141.4
$ws.Range("K33").Value2 = 141.4
$ws.Range("M33").Value2 = 87.59999999999999
$ws.Range("H53").Value2 = 284.16666
$ws.Range("I53").Value2 = 288.8889
$ws.Range("K53").Value2 = 288.8889
$ws.Range("M53").Value2 = 348.1111
$ws.Range("H70").Value2 = 3253.6924
$ws.Range("I70").Value2 = 2025
$ws.Range("J70").Value2 = 3477.0908
$ws.Range("K70").Value2 = 6075
$ws.Range("L70").Value2 = 10431.2724
$ws.Range("M70").Value2 = -5805
$ws.Range("N70").Value2 = -10971.2724
$ws.Range("H73").Value2 = 3253.6924
$ws.Range("I73").Value2 = 2025
$ws.Range("J73").Value2 = 3477.0908
$ws.Range("K73").Value2 = 6075
$ws.Range("L73").Value2 = 10431.2724
$ws.Range("M73").Value2 = -5139
$ws.Range("N73").Value2 = -12303.2724
$ws.Range("H132").Value2 = 968.7037
$ws.Range("I132").Value2 = 890.3461
$ws.Range("K132").Value2 = 2671.0383
$ws.Range("M132").Value2 = -141.0383000000002
$ws.Range("H137").Value2 = 6278
$ws.Range("I137").Value2 = 2307.3809
$ws.Range("K137").Value2 = 6922.1427
$ws.Range("M137").Value2 = -4372.1427
$ws.Range("H138").Value2 = 2853.7646
$ws.Range("I138").Value2 = 2697.1292
$ws.Range("J138").Value2 = 4472.3335
$ws.Range("K138").Value2 = 8091.3876
$ws.Range("L138").Value2 = 13417.0005
$ws.Range("M138").Value2 = -2951.3876
$ws.Range("N138").Value2 = -23697.0005
$ws.Range("H141").Value2 = 10299.8
$ws.Range("I141").Value2 = 3749.5
$ws.Range("K141").Value2 = 11248.5
$ws.Range("M141").Value2 = -6068.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value2 = 1863
$ws.Range("I45").Value2 = 1490.6842
$ws.Range("K45").Value2 = 1490.6842
$ws.Range("M45").Value2 = -1113.6842
$ws.Range("H122").Value2 = 2720
$ws.Range("I122").Value2 = 2080
$ws.Range("J122").Value2 = 2933.3333
$ws.Range("K122").Value2 = 6240
$ws.Range("L122").Value2 = 8799.999899999999
$ws.Range("M122").Value2 = -3790
$ws.Range("N122").Value2 = -13699.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value2 = 4577.1665
$ws.Range("I86").Value2 = 5949.25
$ws.Range("K86").Value2 = 5949.25
$ws.Range("M86").Value2 = -4826.25
$ws.Range("H88").Value2 = 23511.25
$ws.Range("J88").Value2 = 24348.334
$ws.Range("L88").Value2 = 24348.334
$ws.Range("N88").Value2 = -25160.334
$ws.Range("H89").Value2 = 4577.1665
$ws.Range("I89").Value2 = 5949.25
$ws.Range("K89").Value2 = 29746.25
$ws.Range("M89").Value2 = -24130.25
$ws.Range("H91").Value2 = 23511.25
$ws.Range("J91").Value2 = 24348.334
$ws.Range("L91").Value2 = 24348.334
$ws.Range("N91").Value2 = -27156.334
$ws.Range("H94").Value2 = 2764.3713
$ws.Range("J94").Value2 = 3056.125
$ws.Range("L94").Value2 = 3056.125
$ws.Range("N94").Value2 = -3958.125
$ws.Range("H97").Value2 = 26237.666
$ws.Range("I97").Value2 = 14356.5
$ws.Range("K97").Value2 = 14356.5
$ws.Range("M97").Value2 = -13365.5
$ws.Range("H105").Value2 = 6992
$ws.Range("I105").Value2 = 5858
$ws.Range("J105").Value2 = 7685
$ws.Range("K105").Value2 = 5858
$ws.Range("L105").Value2 = 7685
$ws.Range("M105").Value2 = -4111
$ws.Range("N105").Value2 = -11179
$ws.Range("H134").Value2 = 1572.5186
$ws.Range("I134").Value2 = 1520.9546
$ws.Range("K134").Value2 = 4562.8638
$ws.Range("M134").Value2 = -2027.8638

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 7455.5713
$ws.Range("I16").Value2 = 8549.5
$ws.Range("J16").Value2 = 5997
$ws.Range("K16").Value2 = 8549.5
$ws.Range("L16").Value2 = 5997
$ws.Range("M16").Value2 = -8262.5
$ws.Range("N16").Value2 = -6571
$ws.Range("H31").Value2 = 2018.6
$ws.Range("I31").Value2 = 1687.3334
$ws.Range("K31").Value2 = 1687.3334
$ws.Range("M31").Value2 = -1392.3334
$ws.Range("H34").Value2 = 2018.6
$ws.Range("I34").Value2 = 1687.3334
$ws.Range("K34").Value2 = 1687.3334
$ws.Range("M34").Value2 = -1485.3334
$ws.Range("H58").Value2 = 7136.9375
$ws.Range("I58").Value2 = 6320.5557
$ws.Range("K58").Value2 = 6320.5557
$ws.Range("M58").Value2 = -6117.5557
$ws.Range("H86").Value2 = 14724.044
$ws.Range("I86").Value2 = 10933.8
$ws.Range("J86").Value2 = 15776.889
$ws.Range("K86").Value2 = 10933.8
$ws.Range("L86").Value2 = 15776.889
$ws.Range("M86").Value2 = -9810.799999999999
$ws.Range("N86").Value2 = -18022.889
$ws.Range("H89").Value2 = 14724.044
$ws.Range("I89").Value2 = 10933.8
$ws.Range("J89").Value2 = 15776.889
$ws.Range("K89").Value2 = 54669
$ws.Range("L89").Value2 = 78884.44499999999
$ws.Range("M89").Value2 = -49053
$ws.Range("N89").Value2 = -90116.44499999999
$ws.Range("H113").Value2 = 7455.5713
$ws.Range("I113").Value2 = 8549.5
$ws.Range("J113").Value2 = 5997
$ws.Range("K113").Value2 = 8549.5
$ws.Range("L113").Value2 = 5997
$ws.Range("M113").Value2 = -6379.5
$ws.Range("N113").Value2 = -10337
$ws.Range("H134").Value2 = 1937.1765
$ws.Range("I134").Value2 = 1726
$ws.Range("J134").Value2 = 3521
$ws.Range("K134").Value2 = 5178
$ws.Range("L134").Value2 = 10563
$ws.Range("M134").Value2 = -2643
$ws.Range("N134").Value2 = -15633
$ws.Range("H136").Value2 = 7136.9375
$ws.Range("I136").Value2 = 6320.5557
$ws.Range("K136").Value2 = 18961.6671
$ws.Range("M136").Value2 = -16411.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value2 = 3700.25
$ws.Range("J80").Value2 = 3700.25
$ws.Range("L80").Value2 = 11100.75
$ws.Range("N80").Value2 = -12972.75
$ws.Range("H83").Value2 = 3700.25
$ws.Range("J83").Value2 = 3700.25
$ws.Range("L83").Value2 = 33302.25
$ws.Range("N83").Value2 = -42662.25
$ws.Range("H131").Value2 = 6137.857
$ws.Range("I131").Value2 = 1649.8
$ws.Range("K131").Value2 = 4949.4
$ws.Range("M131").Value2 = 90.60000000000036

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I113").Value2 = 1488.8889
$ws.Range("J113").Value2 = 699
$ws.Range("K113").Value2 = 1488.8889
$ws.Range("L113").Value2 = 699
$ws.Range("M113").Value2 = 681.1111000000001
$ws.Range("N113").Value2 = -5039
$ws.Range("H122").Value2 = 2098.9167
$ws.Range("I122").Value2 = 924.625
$ws.Range("J122").Value2 = 2686.0625
$ws.Range("K122").Value2 = 2773.875
$ws.Range("L122").Value2 = 8058.1875
$ws.Range("M122").Value2 = -323.875
$ws.Range("N122").Value2 = -12958.1875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 4099
$ws.Range("I7").Value2 = 3500
$ws.Range("J7").Value2 = 4997.5
$ws.Range("K7").Value2 = 3500
$ws.Range("L7").Value2 = 4997.5
$ws.Range("M7").Value2 = -3388
$ws.Range("N7").Value2 = -5221.5
$ws.Range("H61").Value2 = 10413.214
$ws.Range("I61").Value2 = 12310.111
$ws.Range("K61").Value2 = 12310.111
$ws.Range("M61").Value2 = -12108.111
$ws.Range("H68").Value2 = 6424.6
$ws.Range("I68").Value2 = 5679.933
$ws.Range("J68").Value2 = 7541.6
$ws.Range("K68").Value2 = 5679.933
$ws.Range("L68").Value2 = 7541.6
$ws.Range("M68").Value2 = -4930.933
$ws.Range("N68").Value2 = -9039.6
$ws.Range("H71").Value2 = 6424.6
$ws.Range("I71").Value2 = 5679.933
$ws.Range("J71").Value2 = 7541.6
$ws.Range("K71").Value2 = 28399.665
$ws.Range("L71").Value2 = 37708
$ws.Range("M71").Value2 = -24655.665
$ws.Range("N71").Value2 = -45196
$ws.Range("H93").Value2 = 3730.1538
$ws.Range("I93").Value2 = 2762.5
$ws.Range("K93").Value2 = 2762.5
$ws.Range("M93").Value2 = -1514.5
$ws.Range("H113").Value2 = 10413.214
$ws.Range("I113").Value2 = 12310.111
$ws.Range("K113").Value2 = 12310.111
$ws.Range("M113").Value2 = -10140.111
$ws.Range("H126").Value2 = 4099
$ws.Range("I126").Value2 = 3500
$ws.Range("J126").Value2 = 4997.5
$ws.Range("K126").Value2 = 10500
$ws.Range("L126").Value2 = 14992.5
$ws.Range("M126").Value2 = -8030
$ws.Range("N126").Value2 = -19932.5
$ws.Range("H136").Value2 = 14000
$ws.Range("I136").Value2 = 0
$ws.Range("J136").Value2 = 14000
$ws.Range("K136").Value2 = 0
$ws.Range("L136").Value2 = 42000
$ws.Range("M136").ClearContents() | Out-Null
$ws.Range("N136").Value2 = -47100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value2 = 2390.5881
$ws.Range("I107").Value2 = 1471.9231
$ws.Range("J107").Value2 = 5376.25
$ws.Range("K107").Value2 = 4415.7693
$ws.Range("L107").Value2 = 16128.75
$ws.Range("M107").Value2 = -2495.7693
$ws.Range("N107").Value2 = -19968.75
$ws.Range("H113").Value2 = 2531.724
$ws.Range("I113").Value2 = 2091.842
$ws.Range("K113").Value2 = 6275.526
$ws.Range("M113").Value2 = -4105.526
$ws.Range("H122").Value2 = 35269.44
$ws.Range("I122").Value2 = 1736.96
$ws.Range("J122").Value2 = 128415.22
$ws.Range("K122").Value2 = 5210.88
$ws.Range("L122").Value2 = 385245.66
$ws.Range("M122").Value2 = -2760.88
$ws.Range("N122").Value2 = -390145.66
